$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 7359
$ws.Range("E2").Value = 519
$ws.Range("F2").Value = 519
$ws.Range("G2").Value = 537
$ws.Range("H2").Value = 305
$ws.Range("I2").Value = 315
$ws.Range("J2").Value = -11
$ws.Range("K2").Value = 6943
$ws.Range("L2").Value = 2388
$ws.Range("M2").Value = 4555
$ws.Range("N2").Value = 4431
$ws.Range("O2").Value = 124
$ws.Range("P2").Value = 290
$ws.Range("Q2").Value = 259
$ws.Range("R2").Value = -1294
$ws.Range("S2").Value = 995
$ws.Range("T2").Value = 460
$ws.Range("U2").Value = -201
$ws.Range("V2").Value = 1444
$ws.Range("W2").Value = 7.06
$ws.Range("X2").Value = 4.14
$ws.Range("Y2").Value = 7.2
$ws.Range("Z2").Value = 4.79
$ws.Range("AA2").Value = 52.42
$ws.Range("AB2").Value = 1573.33
$ws.Range("AC2").Value = 2722
$ws.Range("AD2").Value = 24.68
$ws.Range("AE2").Value = 42930
$ws.Range("AF2").Value = 1.57
$ws.Range("AG2").Value = 700
$ws.Range("AH2").Value = 1.04
$ws.Range("AI2").Value = 22.9
$ws.Range("AJ2").Value = 11586575

# Row 3 updates
$ws.Range("D3").Value = 8397
$ws.Range("E3").Value = 436
$ws.Range("F3").Value = 436
$ws.Range("G3").Value = 478
$ws.Range("H3").Value = 357
$ws.Range("I3").Value = 429
$ws.Range("J3").Value = -73
$ws.Range("K3").Value = 9445
$ws.Range("L3").Value = 3643
$ws.Range("M3").Value = 5802
$ws.Range("N3").Value = 4801
$ws.Range("O3").Value = 1002
$ws.Range("P3").Value = 290
$ws.Range("Q3").Value = 359
$ws.Range("R3").Value = -1086
$ws.Range("S3").Value = 679
$ws.Range("T3").Value = 951
$ws.Range("U3").Value = -592
$ws.Range("V3").Value = 2369
$ws.Range("W3").Value = 5.19
$ws.Range("X3").Value = 4.25
$ws.Range("Y3").Value = 9.3
$ws.Range("Z3").Value = 4.35
$ws.Range("AA3").Value = 62.78
$ws.Range("AB3").Value = 1708.12
$ws.Range("AC3").Value = 3705
$ws.Range("AD3").Value = 19.19
$ws.Range("AE3").Value = 46515
$ws.Range("AF3").Value = 1.53
$ws.Range("AG3").Value = 700
$ws.Range("AH3").Value = 0.98
$ws.Range("AI3").Value = 16.83
$ws.Range("AJ3").Value = 11586575

# Row 4 updates
$ws.Range("D4").Value = 8839
$ws.Range("E4").Value = 259
$ws.Range("F4").Value = 259
$ws.Range("G4").Value = 276
$ws.Range("H4").Value = 261
$ws.Range("I4").Value = 270
$ws.Range("J4").Value = -8
$ws.Range("K4").Value = 11447
$ws.Range("L4").Value = 5401
$ws.Range("M4").Value = 6046
$ws.Range("N4").Value = 4995
$ws.Range("O4").Value = 1051
$ws.Range("P4").Value = 290
$ws.Range("Q4").Value = 154
$ws.Range("R4").Value = -1558
$ws.Range("S4").Value = 1592
$ws.Range("T4").Value = 1414
$ws.Range("U4").Value = -1260
$ws.Range("V4").Value = 3969
$ws.Range("W4").Value = 2.93
$ws.Range("X4").Value = 2.96
$ws.Range("Y4").Value = 5.51
$ws.Range("Z4").Value = 2.5
$ws.Range("AA4").Value = 89.35
$ws.Range("AB4").Value = 1781.67
$ws.Range("AC4").Value = 2328
$ws.Range("AD4").Value = 29.51
$ws.Range("AE4").Value = 48396
$ws.Range("AF4").Value = 1.42
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 0.87
$ws.Range("AI4").Value = 22.96
$ws.Range("AJ4").Value = 11586575

# Row 5 updates
$ws.Range("D5").Value = 9603
$ws.Range("E5").Value = 390
$ws.Range("F5").Value = 390
$ws.Range("G5").Value = 361
$ws.Range("H5").Value = 354
$ws.Range("I5").Value = 316
$ws.Range("J5").Value = 39
$ws.Range("K5").Value = 12135
$ws.Range("L5").Value = 5820
$ws.Range("M5").Value = 6315
$ws.Range("N5").Value = 5194
$ws.Range("O5").Value = 1121
$ws.Range("P5").Value = 290
$ws.Range("Q5").Value = 1039
$ws.Range("R5").Value = -730
$ws.Range("S5").Value = -167
$ws.Range("T5").Value = 691
$ws.Range("U5").Value = 349
$ws.Range("V5").Value = 3851
$ws.Range("W5").Value = 4.06
$ws.Range("X5").Value = 3.69
$ws.Range("Y5").Value = 6.2
$ws.Range("Z5").Value = 3
$ws.Range("AA5").Value = 92.15
$ws.Range("AB5").Value = 1859.3
$ws.Range("AC5").Value = 2726
$ws.Range("AD5").Value = 60.35
$ws.Range("AE5").Value = 50329
$ws.Range("AF5").Value = 3.27
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 0.36
$ws.Range("AI5").Value = 19.61
$ws.Range("AJ5").Value = 11586575

# Row 6 updates
$ws.Range("D6").Value = 10314
$ws.Range("E6").Value = 246
$ws.Range("F6").Value = 246
$ws.Range("G6").Value = -6
$ws.Range("H6").Value = -154
$ws.Range("I6").Value = -152
$ws.Range("K6").Value = 12078
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 6078
$ws.Range("N6").Value = 4941
$ws.Range("P6").Value = 290
$ws.Range("Q6").Value = 475
$ws.Range("R6").Value = -686
$ws.Range("S6").Value = -48
$ws.Range("T6").Value = 570
$ws.Range("U6").Value = -95
$ws.Range("V6").Value = 3868
$ws.Range("W6").Value = 2.38
$ws.Range("X6").Value = -1.49
$ws.Range("Y6").Value = -2.99
$ws.Range("Z6").Value = -1.27
$ws.Range("AA6").Value = 98.72
$ws.Range("AB6").Value = 1793.26
$ws.Range("AC6").Value = -1308
$ws.Range("AD6").Value = -144.15
$ws.Range("AE6").Value = 47877
$ws.Range("AF6").Value = 3.94
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 0.32
$ws.Range("AI6").Value = -40.87
$ws.Range("AJ6").Value = 11586575

# Row 7: clear removed cells, update remaining
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("AG7").Value = 600
$ws.Range("AH7").Value = 0.51

# Row 8: clear removed cells, update remaining
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("AG8").Value = 600
$ws.Range("AH8").Value = 0.51

# Row 9: clear removed cells, update remaining
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AI9").ClearContents()
$ws.Range("AG9").Value = 562
$ws.Range("AH9").Value = 0.48
